# Auto-generated edit script: updates horarios (schedule) data for Linea 141
# Reflects a refreshed scrape: new "Ultima actualizacion" timestamp, updated
# rows, and 8 newly appended arrival rows (265-272) on sheet "LP1912".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 16:44:58"

$ws1.Range("A3").Value = "Total filas: 267"

$ws1.Range("C38").Value = "11_ETCHEVERRY"

$ws1.Range("C39").Value = "15_ABASTO"

$ws1.Range("C106").Value = "16_SANTA ANA"

$ws1.Range("C107").Value = "10_OLMOS"

$ws1.Range("A139").Value = "11:53:44"
$ws1.Range("C139").Value = "23_HERNANDEZ"
$ws1.Range("D139").Value = 43

$ws1.Range("A140").Value = "10:49:38"
$ws1.Range("C140").Value = "27_EL RETIRO"
$ws1.Range("D140").Value = 107

$ws1.Range("C204").Value = "15X38_ABASTO"

$ws1.Range("A205").Value = "14:32:44"
$ws1.Range("C205").Value = "10_OLMOS"
$ws1.Range("D205").Value = 81

$ws1.Range("A206").Value = "13:55:43"
$ws1.Range("C206").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D206").Value = 118

$ws1.Range("A236").Value = "16:44:58"
$ws1.Range("B236").Value = "17:06"
$ws1.Range("D236").Value = 22

$ws1.Range("A237").Value = "15:16:46"
$ws1.Range("B237").Value = "17:07"
$ws1.Range("C237").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D237").Value = 111

$ws1.Range("A238").Value = "16:28:21"
$ws1.Range("B238").Value = "17:08"
$ws1.Range("C238").Value = "10_OLMOS"
$ws1.Range("D238").Value = 40

$ws1.Range("A239").Value = "15:16:46"
$ws1.Range("B239").Value = "17:09"
$ws1.Range("D239").Value = 113

$ws1.Range("A240").Value = "16:12:06"
$ws1.Range("B240").Value = "17:10"
$ws1.Range("C240").Value = "215C_EL PATO"
$ws1.Range("D240").Value = 58

$ws1.Range("A241").Value = "16:44:58"
$ws1.Range("B241").Value = "17:20"
$ws1.Range("C241").Value = "15X38_ABASTO"
$ws1.Range("D241").Value = 36

$ws1.Range("A242").Value = "15:44:42"
$ws1.Range("B242").Value = "17:21"
$ws1.Range("C242").Value = "15X38_ABASTO"

$ws1.Range("A243").Value = "16:28:21"
$ws1.Range("B243").Value = "17:32"
$ws1.Range("C243").Value = "27_EL RETIRO"
$ws1.Range("D243").Value = 64

$ws1.Range("A244").Value = "15:56:56"
$ws1.Range("B244").Value = "17:33"
$ws1.Range("C244").Value = "17_ROMERO"
$ws1.Range("D244").Value = 97

$ws1.Range("A245").Value = "15:44:42"
$ws1.Range("B245").Value = "17:34"
$ws1.Range("C245").Value = "17_ROMERO"
$ws1.Range("D245").Value = 110

$ws1.Range("B246").Value = "17:36"
$ws1.Range("C246").Value = "27_EL RETIRO"
$ws1.Range("D246").Value = 112

$ws1.Range("B247").Value = "17:37"
$ws1.Range("C247").Value = "27_EL RETIRO"
$ws1.Range("D247").Value = 85

$ws1.Range("A248").Value = "15:44:42"
$ws1.Range("B248").Value = "17:38"
$ws1.Range("C248").Value = "215B_EL PATO"
$ws1.Range("D248").Value = 114

$ws1.Range("B249").Value = "17:39"
$ws1.Range("C249").Value = "215B_EL PATO"
$ws1.Range("D249").Value = 87

$ws1.Range("A250").Value = "15:56:56"
$ws1.Range("B250").Value = "17:45"
$ws1.Range("C250").Value = "215_EL PELIGRO"
$ws1.Range("D250").Value = 109

$ws1.Range("A251").Value = "16:12:06"
$ws1.Range("B251").Value = "17:46"
$ws1.Range("C251").Value = "215_EL PELIGRO"
$ws1.Range("D251").Value = 94

$ws1.Range("A252").Value = "16:12:06"
$ws1.Range("B252").Value = "17:49"
$ws1.Range("C252").Value = "10_OLMOS"
$ws1.Range("D252").Value = 97

$ws1.Range("A253").Value = "15:56:56"
$ws1.Range("B253").Value = "17:51"
$ws1.Range("C253").Value = "10_OLMOS"
$ws1.Range("D253").Value = 115

$ws1.Range("A254").Value = "16:28:21"
$ws1.Range("B254").Value = "17:52"
$ws1.Range("C254").Value = "23_HERNANDEZ"
$ws1.Range("D254").Value = 84

$ws1.Range("A255").Value = "16:37:37"
$ws1.Range("B255").Value = "17:53"
$ws1.Range("C255").Value = "23_HERNANDEZ"
$ws1.Range("D255").Value = 76

$ws1.Range("A256").Value = "16:44:58"
$ws1.Range("B256").Value = "17:57"
$ws1.Range("C256").Value = "17_ROMERO"
$ws1.Range("D256").Value = 73

$ws1.Range("B257").Value = "17:58"
$ws1.Range("C257").Value = "17_ROMERO"
$ws1.Range("D257").Value = 106

$ws1.Range("A258").Value = "16:28:21"
$ws1.Range("B258").Value = "18:05"
$ws1.Range("C258").Value = "11_ETCHEVERRY"
$ws1.Range("D258").Value = 97

$ws1.Range("A259").Value = "16:12:06"
$ws1.Range("B259").Value = "18:06"
$ws1.Range("C259").Value = "11_ETCHEVERRY"
$ws1.Range("D259").Value = 114

$ws1.Range("A260").Value = "16:44:58"
$ws1.Range("B260").Value = "18:09"
$ws1.Range("C260").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D260").Value = 85

$ws1.Range("A261").Value = "16:44:58"
$ws1.Range("B261").Value = "18:09"
$ws1.Range("C261").Value = "15_ABASTO"
$ws1.Range("D261").Value = 85

$ws1.Range("A262").Value = "16:12:06"
$ws1.Range("B262").Value = "18:10"
$ws1.Range("C262").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D262").Value = 118

$ws1.Range("A263").Value = "16:12:06"
$ws1.Range("B263").Value = "18:10"
$ws1.Range("C263").Value = "15_ABASTO"
$ws1.Range("D263").Value = 118

$ws1.Range("A264").Value = "16:44:58"
$ws1.Range("B264").Value = "18:16"
$ws1.Range("C264").Value = "10_OLMOS"
$ws1.Range("D264").Value = 92

$ws1.Range("A265").Value = "16:28:21"
$ws1.Range("B265").Value = "18:17"
$ws1.Range("C265").Value = "10_OLMOS"
$ws1.Range("D265").Value = 109
$ws1.Range("E265").Value = "LP1912"

$ws1.Range("A266").Value = "16:37:37"
$ws1.Range("B266").Value = "18:21"
$ws1.Range("C266").Value = "215C_EL PATO"
$ws1.Range("D266").Value = 104
$ws1.Range("E266").Value = "LP1912"

$ws1.Range("A267").Value = "16:28:21"
$ws1.Range("B267").Value = "18:22"
$ws1.Range("C267").Value = "215C_EL PATO"
$ws1.Range("D267").Value = 114
$ws1.Range("E267").Value = "LP1912"

$ws1.Range("A268").Value = "16:28:21"
$ws1.Range("B268").Value = "18:25"
$ws1.Range("C268").Value = "16_SANTA ANA"
$ws1.Range("D268").Value = 117
$ws1.Range("E268").Value = "LP1912"

$ws1.Range("A269").Value = "16:37:37"
$ws1.Range("B269").Value = "18:29"
$ws1.Range("C269").Value = "14_ABASTO"
$ws1.Range("D269").Value = 112
$ws1.Range("E269").Value = "LP1912"

$ws1.Range("A270").Value = "16:44:58"
$ws1.Range("B270").Value = "18:35"
$ws1.Range("C270").Value = "15X38_ABASTO"
$ws1.Range("D270").Value = 111
$ws1.Range("E270").Value = "LP1912"

$ws1.Range("A271").Value = "16:37:37"
$ws1.Range("B271").Value = "18:36"
$ws1.Range("C271").Value = "15X38_ABASTO"
$ws1.Range("D271").Value = 119
$ws1.Range("E271").Value = "LP1912"

$ws1.Range("A272").Value = "16:44:58"
$ws1.Range("B272").Value = "18:40"
$ws1.Range("C272").Value = "10_OLMOS"
$ws1.Range("D272").Value = 116
$ws1.Range("E272").Value = "LP1912"

# Sheets "LP1912-215" and "6203-6173" only need their "Ultima actualizacion" stamp updated
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 16:44:58"

$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 16:44:58"

